# Remove the "Ι.Κ.Υ." list-item paragraph from the ΚΟΙΝΟΠΟΙΗΣΗ (notification)
# list. Deleting the paragraph's Range (which, for a Paragraph object,
# includes the trailing paragraph mark) merges it with the following
# paragraph, exactly as shown in the diff.

$d = $word.ActiveDocument

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "Ι\.Κ\.Υ\.") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $target.Range.Delete()
}
